$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = "[53.39136156740169, 71.79568969667945]"
$ws.Range("U2").Value = "[42.58496710922837, 54.955930866258896]"

# Row 3
$ws.Range("M3").Value = "[50.68759602986137, 74.84054149884791]"
$ws.Range("N3").Value = [double]"1.216804434989172e-13"
$ws.Range("O3").Value = [double]"1.216804434989172e-13"
$ws.Range("U3").Value = "[43.65908750583876, 57.03938032983465]"

# Row 4
$ws.Range("M4").Value = "[51.588456789048955, 77.0827936739106]"
$ws.Range("N4").Value = [double]"3.104183576851938e-13"
$ws.Range("O4").Value = [double]"3.104183576851938e-13"
$ws.Range("U4").Value = "[44.796859958799374, 58.276188327118604]"

# Row 5
$ws.Range("M5").Value = "[53.461447658397226, 77.52886552882813]"
$ws.Range("N5").Value = [double]"2.708944180085382e-14"
$ws.Range("O5").Value = [double]"2.708944180085382e-14"
$ws.Range("Q5").Value = "[0.1823947686768861, 0.5597632555945768]"
$ws.Range("R5").Value = [double]"0.0002632491671181736"
$ws.Range("S5").Value = [double]"0.0002632491671181736"
$ws.Range("U5").Value = "[43.438322371081114, 56.79538285293016]"
$ws.Range("Y5").Value = [double]"23.31021021021073"
$ws.Range("Z5").Value = [double]"24.84714714714769"

# Row 6
$ws.Range("M6").Value = "[55.043738528412234, 75.45552518149921]"
$ws.Range("N6").Value = [double]"2.220446049250313e-16"
$ws.Range("O6").Value = [double]"2.220446049250313e-16"
$ws.Range("U6").Value = "[43.30722592905357, 56.70940822191007]"

# Row 7
$ws.Range("M7").Value = "[51.68065666652505, 75.98027726608056]"
$ws.Range("N7").Value = [double]"8.593126210598712e-14"
$ws.Range("O7").Value = [double]"8.593126210598712e-14"
$ws.Range("U7").Value = "[42.52421094244815, 56.0196110186253]"

# Row 8
$ws.Range("M8").Value = "[50.24206945724285, 74.78436814004418]"
$ws.Range("N8").Value = [double]"2.309263891220326e-13"
$ws.Range("O8").Value = [double]"2.309263891220326e-13"
$ws.Range("U8").Value = "[42.782203231339906, 55.61685508291853]"

# Row 9
$ws.Range("M9").Value = "[53.129106754591604, 73.18538963112829]"
$ws.Range("N9").Value = [double]"2.220446049250313e-16"
$ws.Range("O9").Value = [double]"2.220446049250313e-16"
$ws.Range("U9").Value = "[42.26247941425813, 54.07459049475853]"

# Row 10
$ws.Range("M10").Value = "[53.50451482907803, 71.89482636908372]"
$ws.Range("U10").Value = "[45.10268176544007, 57.511254288606295]"

# Row 11
$ws.Range("M11").Value = "[51.346574514509115, 74.31095993515582]"
$ws.Range("N11").Value = [double]"2.264854970235319e-14"
$ws.Range("O11").Value = [double]"2.264854970235319e-14"
$ws.Range("U11").Value = "[43.20491024349958, 56.53887521885116]"

# Row 12
$ws.Range("M12").Value = "[51.53880904866267, 76.95503358623168]"
$ws.Range("N12").Value = [double]"2.942091015256665e-13"
$ws.Range("O12").Value = [double]"2.942091015256665e-13"
$ws.Range("U12").Value = "[42.07585227387733, 55.50616214103171]"

# Row 13
$ws.Range("M13").Value = "[53.52520576219777, 77.37610577351452]"
$ws.Range("N13").Value = [double]"2.042810365310288e-14"
$ws.Range("O13").Value = [double]"2.042810365310288e-14"
$ws.Range("U13").Value = "[41.55667895965422, 54.86585129873916]"

# Row 14
$ws.Range("M14").Value = "[55.045386459209865, 75.85596276865168]"
$ws.Range("U14").Value = "[41.55625355381809, 54.86670304501221]"
